$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10: "Console write line / print" (adds a shared string + sheet row)
$ws.Range("A10").Value = "Console write line / print"

# New cell comment on A10, matching the other "Jeremy:" style comments already
# on the sheet (bold author prefix + code snippet).
$commentText = "Jeremy:" + [char]10 + 'System.out.println("FizzBuzz");'
$comment = $ws.Range("A10").AddComment($commentText)

# Match the author's final selection (cell T11) recorded in the sheet view.
$ws.Range("T11").Select() | Out-Null
